$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 88
$ws.Range("H88").Value = 1983.091
$ws.Range("I88").Value = 1849.75
$ws.Range("J88").Value = 2059.2856
$ws.Range("K88").Value = 1849.75
$ws.Range("L88").Value = 2059.2856
$ws.Range("M88").Value = -1443.75
$ws.Range("N88").Value = -2871.2856
# Row 91
$ws.Range("H91").Value = 1983.091
$ws.Range("I91").Value = 1849.75
$ws.Range("J91").Value = 2059.2856
$ws.Range("K91").Value = 1849.75
$ws.Range("L91").Value = 2059.2856
$ws.Range("M91").Value = -445.75
$ws.Range("N91").Value = -4867.2856
# Row 111
$ws.Range("H111").Value = 2707.25
# Row 129
$ws.Range("H129").Value = 218196.48
$ws.Range("I129").Value = 339.4
$ws.Range("J129").Value = 244764.42
$ws.Range("K129").Value = 1018.2
$ws.Range("L129").Value = 734293.26
$ws.Range("M129").Value = 3981.8
$ws.Range("N129").Value = -744293.26
# Row 132
$ws.Range("H132").Value = 2088.9556
$ws.Range("I132").Value = 2090.9773
$ws.Range("J132").Value = 2000
$ws.Range("K132").Value = 6272.9319
$ws.Range("L132").Value = 6000
$ws.Range("M132").Value = -3742.9319
# Row 137
$ws.Range("H137").Value = 2897.842
$ws.Range("I137").Value = 2797.4375
$ws.Range("J137").Value = 3433.3333
$ws.Range("K137").Value = 8392.3125
$ws.Range("L137").Value = 10299.9999
$ws.Range("M137").Value = -5842.3125
$ws.Range("N137").Value = -15399.9999

$ws = $wb.Worksheets.Item("ARM")
# Row 32
$ws.Range("H32").Value = 4158.949
$ws.Range("I32").Value = 4149.357
$ws.Range("J32").Value = 4338
$ws.Range("K32").Value = 4149.357
$ws.Range("L32").Value = 4338
$ws.Range("M32").Value = -3862.357
$ws.Range("N32").Value = -4912
# Row 63
$ws.Range("H63").Value = 2124.75
$ws.Range("I63").Value = 1999.6666
$ws.Range("J63").Value = 2500
$ws.Range("K63").Value = 1999.6666
$ws.Range("L63").Value = 2500
$ws.Range("M63").Value = -1313.6666
# Row 66
$ws.Range("H66").Value = 2124.75
$ws.Range("I66").Value = 1999.6666
$ws.Range("J66").Value = 2500
$ws.Range("K66").Value = 9998.333000000001
$ws.Range("L66").Value = 12500
$ws.Range("M66").Value = -6566.333000000001
# Row 74
$ws.Range("H74").Value = 66667972
$ws.Range("I74").Value = 100000710
$ws.Range("J74").Value = 2500
$ws.Range("K74").Value = 100000710
$ws.Range("L74").Value = 2500
$ws.Range("M74").Value = -99999836
$ws.Range("N74").Value = -4248
# Row 77
$ws.Range("H77").Value = 66667972
$ws.Range("I77").Value = 100000710
$ws.Range("J77").Value = 2500
$ws.Range("K77").Value = 500003550
$ws.Range("L77").Value = 12500
$ws.Range("M77").Value = -499999182
$ws.Range("N77").Value = -21236
# Row 109
$ws.Range("H109").Value = 20000
$ws.Range("I109").Value = 0
$ws.Range("J109").Value = 20000
$ws.Range("K109").Value = 0
$ws.Range("L109").Value = 20000
$ws.Range("N109").Value = -22774
# Row 112
$ws.Range("H112").Value = 37000
$ws.Range("I112").Value = 0
$ws.Range("J112").Value = 37000
$ws.Range("K112").Value = 0
$ws.Range("L112").Value = 37000
$ws.Range("N112").Value = -39954
# Row 132
$ws.Range("H132").Value = 12815.422
$ws.Range("I132").Value = 1640.4849
$ws.Range("J132").Value = 43546.5
$ws.Range("K132").Value = 4921.4547
$ws.Range("L132").Value = 130639.5
$ws.Range("M132").Value = -2391.4547
$ws.Range("N132").Value = -135699.5

$ws = $wb.Worksheets.Item("CRP")
# Row 16
$ws.Range("H16").Value = 885.55554
$ws.Range("I16").Value = 852.8570999999999
$ws.Range("J16").Value = 1000
$ws.Range("K16").Value = 852.8570999999999
$ws.Range("L16").Value = 1000
$ws.Range("M16").Value = -565.8570999999999
# Row 31
$ws.Range("H31").Value = 12680.366
$ws.Range("I31").Value = 21505.05
$ws.Range("J31").Value = 4275.905
$ws.Range("K31").Value = 21505.05
$ws.Range("L31").Value = 4275.905
$ws.Range("M31").Value = -21210.05
$ws.Range("N31").Value = -4865.905
# Row 34
$ws.Range("H34").Value = 12680.366
$ws.Range("I34").Value = 21505.05
$ws.Range("J34").Value = 4275.905
$ws.Range("K34").Value = 21505.05
$ws.Range("L34").Value = 4275.905
$ws.Range("M34").Value = -21303.05
$ws.Range("N34").Value = -4679.905
# Row 51
$ws.Range("H51").Value = 14990
$ws.Range("I51").Value = 0
$ws.Range("J51").Value = 14990
$ws.Range("K51").Value = 0
$ws.Range("L51").Value = 14990
$ws.Range("N51").Value = -16462
$ws.Range("M51").ClearContents()
# Row 61
$ws.Range("H61").Value = 14990
$ws.Range("I61").Value = 0
$ws.Range("J61").Value = 14990
$ws.Range("K61").Value = 0
$ws.Range("L61").Value = 14990
$ws.Range("N61").Value = -15686
$ws.Range("M61").ClearContents()
# Row 86
$ws.Range("H86").Value = 9269656
$ws.Range("I86").Value = 2479
$ws.Range("J86").Value = 23832362
$ws.Range("K86").Value = 2479
$ws.Range("L86").Value = 23832362
$ws.Range("M86").Value = -1356
# Row 89
$ws.Range("H89").Value = 9269656
$ws.Range("I89").Value = 2479
$ws.Range("J89").Value = 23832362
$ws.Range("K89").Value = 12395
$ws.Range("L89").Value = 119161810
$ws.Range("M89").Value = -6779
# Row 113
$ws.Range("H113").Value = 885.55554
$ws.Range("I113").Value = 852.8570999999999
$ws.Range("J113").Value = 1000
$ws.Range("K113").Value = 852.8570999999999
$ws.Range("L113").Value = 1000
$ws.Range("M113").Value = 1317.1429
# Row 132
$ws.Range("H132").Value = 11708.942
$ws.Range("I132").Value = 15058.811
$ws.Range("J132").Value = 3445.9333
$ws.Range("K132").Value = 45176.433
$ws.Range("L132").Value = 10337.7999
$ws.Range("M132").Value = -42646.433
$ws.Range("N132").Value = -15397.7999
# Row 134
$ws.Range("H134").Value = 1098.6031
$ws.Range("I134").Value = 837.45
$ws.Range("J134").Value = 1552.7826
$ws.Range("K134").Value = 2512.35
$ws.Range("L134").Value = 4658.3478
$ws.Range("M134").Value = 22.64999999999964

$ws = $wb.Worksheets.Item("CUL")
# Row 2
$ws.Range("H2").Value = 6321.0625
$ws.Range("I2").Value = 14328.571
$ws.Range("J2").Value = 93
$ws.Range("K2").Value = 85971.42600000001
$ws.Range("L2").Value = 558
$ws.Range("M2").Value = -85858.42600000001
$ws.Range("N2").Value = -784
# Row 26
$ws.Range("H26").Value = 366.42856
$ws.Range("I26").Value = 141.25
$ws.Range("J26").Value = 666.6667
$ws.Range("K26").Value = 423.75
$ws.Range("L26").Value = 2000.0001
$ws.Range("M26").Value = -135.75
$ws.Range("N26").Value = -2576.0001
# Row 117
$ws.Range("H117").Value = 2623.6924
$ws.Range("I117").Value = 1119.3334
$ws.Range("J117").Value = 3075
$ws.Range("K117").Value = 3358.0002
$ws.Range("L117").Value = 9225
$ws.Range("M117").Value = 83.99980000000005
$ws.Range("N117").Value = -16109
# Row 131
$ws.Range("H131").Value = 768.89
$ws.Range("I131").Value = 0
$ws.Range("J131").Value = 768.89
$ws.Range("K131").Value = 0
$ws.Range("L131").Value = 2306.67
$ws.Range("N131").Value = -12386.67

$ws = $wb.Worksheets.Item("GSM")
# Row 132
$ws.Range("H132").Value = 20357.768
$ws.Range("I132").Value = 3689.348
$ws.Range("J132").Value = 75125.42999999999
$ws.Range("K132").Value = 11068.044
$ws.Range("L132").Value = 225376.29
$ws.Range("M132").Value = -8538.044
$ws.Range("N132").Value = -230436.29

$ws = $wb.Worksheets.Item("LTW")
# Row 46
$ws.Range("H46").Value = 753.1539
$ws.Range("I46").Value = 636.375
$ws.Range("J46").Value = 940
$ws.Range("K46").Value = 636.375
$ws.Range("L46").Value = 940
$ws.Range("M46").Value = -448.375
$ws.Range("N46").Value = -1316
# Row 50
$ws.Range("H50").Value = 19000
$ws.Range("I50").Value = 0
$ws.Range("J50").Value = 19000
$ws.Range("K50").Value = 0
$ws.Range("L50").Value = 19000
$ws.Range("N50").Value = -20274
# Row 61
$ws.Range("H61").Value = 3532.6553
$ws.Range("I61").Value = 1954.7273
$ws.Range("J61").Value = 8491.857
$ws.Range("K61").Value = 1954.7273
$ws.Range("L61").Value = 8491.857
$ws.Range("M61").Value = -1752.7273
$ws.Range("N61").Value = -8895.857
# Row 110
$ws.Range("H110").Value = 38000
$ws.Range("I110").Value = 0
$ws.Range("J110").Value = 38000
$ws.Range("K110").Value = 0
$ws.Range("L110").Value = 38000
$ws.Range("N110").Value = -46180
# Row 113
$ws.Range("H113").Value = 3532.6553
$ws.Range("I113").Value = 1954.7273
$ws.Range("J113").Value = 8491.857
$ws.Range("K113").Value = 1954.7273
$ws.Range("L113").Value = 8491.857
$ws.Range("M113").Value = 215.2727
$ws.Range("N113").Value = -12831.857

$ws = $wb.Worksheets.Item("WVR")
# Row 81
$ws.Range("H81").Value = 62501268
$ws.Range("I81").Value = 1406.6923
$ws.Range("J81").Value = 333334000
$ws.Range("K81").Value = 2813.3846
$ws.Range("L81").Value = 666668000
$ws.Range("M81").Value = -1752.3846
$ws.Range("N81").Value = -666670122
# Row 84
$ws.Range("H84").Value = 62501268
$ws.Range("I84").Value = 1406.6923
$ws.Range("J84").Value = 333334000
$ws.Range("K84").Value = 14066.923
$ws.Range("L84").Value = 3333340000
$ws.Range("M84").Value = -8762.922999999999
$ws.Range("N84").Value = -3333350608
# Row 100
$ws.Range("H100").Value = 400
$ws.Range("I100").Value = 400
$ws.Range("J100").Value = 0
$ws.Range("K100").Value = 800
$ws.Range("L100").Value = 0
$ws.Range("M100").Value = -259
$ws.Range("N100").ClearContents()
